$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "Big StockRoom" header typo -> "Big Stockroom"
$ws.Range("C1").Value = "Big Stockroom"

# Update stock-done entries for row 2 (Penunt Oil) and row 3 (glass)
$ws.Range("A2").Value = "Penunt Oil"
$ws.Range("B2").Value = 2021
$ws.Range("B3").Value = 1020

# Remove the extra "Rahul" rows (4 and 5) that are no longer needed
$ws.Rows("4:5").Delete()
